$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31498

$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -97488

$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("N131").ClearContents()

$ws.Range("H132").Value = 3857.5
$ws.Range("I132").Value = 4324.4736
$ws.Range("K132").Value = 12973.4208
$ws.Range("M132").Value = -10443.4208

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H135").Value = 18524298
$ws.Range("I135").Value = 592.0952
$ws.Range("J135").Value = 83357260
$ws.Range("K135").Value = 5328.8568
$ws.Range("L135").Value = 750215340
$ws.Range("M135").Value = -2793.8568
$ws.Range("N135").Value = -750220410

$ws.Range("H136").Value = 49780
$ws.Range("J136").Value = 49780
$ws.Range("L136").Value = 49780
$ws.Range("N136").Value = -59980

$ws.Range("H138").Value = 2107.575
$ws.Range("I138").Value = 1732.5555
$ws.Range("J138").Value = 2216.4517
$ws.Range("K138").Value = 5197.666499999999
$ws.Range("L138").Value = 6649.355100000001
$ws.Range("M138").Value = -57.66649999999936
$ws.Range("N138").Value = -16929.3551

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

$ws.Range("H56").Value = 16000
$ws.Range("J56").Value = 16000
$ws.Range("L56").Value = 16000
$ws.Range("N56").Value = -17484

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 13950
$ws.Range("I39").Value = 13950
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 13950
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -13561
$ws.Range("N39").ClearContents()

$ws.Range("H54").Value = 12176.6
$ws.Range("I54").Value = 4441.5
$ws.Range("J54").Value = 17333.334
$ws.Range("K54").Value = 4441.5
$ws.Range("L54").Value = 17333.334
$ws.Range("M54").Value = -3957.5
$ws.Range("N54").Value = -18301.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 134.45454
$ws.Range("I19").Value = 134.45454
$ws.Range("K19").Value = 134.45454
$ws.Range("M19").Value = 35.54545999999999

$ws.Range("H24").Value = 134.45454
$ws.Range("I24").Value = 134.45454
$ws.Range("K24").Value = 134.45454
$ws.Range("M24").Value = 35.54545999999999

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H134").Value = 737.43335
$ws.Range("I134").Value = 605.72
$ws.Range("J134").Value = 1396
$ws.Range("K134").Value = 1817.16
$ws.Range("L134").Value = 4188
$ws.Range("M134").Value = 717.8399999999999
$ws.Range("N134").Value = -9258

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 424.8889
$ws.Range("I10").Value = 424.8889
$ws.Range("K10").Value = 1274.6667
$ws.Range("M10").Value = -1135.6667

$ws.Range("H50").Value = 230
$ws.Range("I50").Value = 325
$ws.Range("J50").Value = 166.66667
$ws.Range("K50").Value = 975
$ws.Range("L50").Value = 500.00001
$ws.Range("M50").Value = -494
$ws.Range("N50").Value = -1462.00001

$ws.Range("H53").Value = 230
$ws.Range("I53").Value = 325
$ws.Range("J53").Value = 166.66667
$ws.Range("K53").Value = 975
$ws.Range("L53").Value = 500.00001
$ws.Range("M53").Value = -494
$ws.Range("N53").Value = -1462.00001

$ws.Range("H122").Value = 685.96
$ws.Range("J122").Value = 867.3125
$ws.Range("L122").Value = 7805.8125
$ws.Range("N122").Value = -12705.8125

$ws.Range("H131").Value = 818
$ws.Range("J131").Value = 818.1818
$ws.Range("L131").Value = 2454.5454
$ws.Range("N131").Value = -12534.5454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4857143
$ws.Range("J7").Value = 4600000
$ws.Range("L7").Value = 4600000
$ws.Range("N7").Value = -4600224

$ws.Range("H8").Value = 4857143
$ws.Range("J8").Value = 4600000
$ws.Range("L8").Value = 4600000
$ws.Range("N8").Value = -4600278

$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws.Range("H12").Value = 6182353
$ws.Range("J12").Value = 3100000
$ws.Range("L12").Value = 3100000
$ws.Range("N12").Value = -3100280

$ws.Range("H17").Value = 850
$ws.Range("J17").Value = 850
$ws.Range("L17").Value = 850
$ws.Range("N17").Value = -1186

$ws.Range("H28").Value = 2496.5
$ws.Range("I28").Value = 13
$ws.Range("J28").Value = 4980
$ws.Range("K28").Value = 13
$ws.Range("L28").Value = 4980
$ws.Range("M28").Value = 179
$ws.Range("N28").Value = -5364

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 1000
$ws.Range("J20").Value = 1000
$ws.Range("L20").Value = 1000
$ws.Range("N20").Value = -1452

$ws.Range("H136").Value = 37689.285
$ws.Range("J136").Value = 3021
$ws.Range("L136").Value = 9063
$ws.Range("N136").Value = -14163

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1812.8572
$ws.Range("I6").Value = 1005
$ws.Range("K6").Value = 1005
$ws.Range("M6").Value = -890

$ws.Range("H96").Value = 1019.2
$ws.Range("J96").Value = 1074
$ws.Range("L96").Value = 1074
$ws.Range("N96").Value = -3820

$ws.Range("H136").Value = 45456670
$ws.Range("I136").Value = 58825616
$ws.Range("K136").Value = 176476848
$ws.Range("M136").Value = -176474298

$ws.Range("H137").Value = 46518.125
$ws.Range("J137").Value = 47429
$ws.Range("L137").Value = 47429
$ws.Range("N137").Value = -57629
